$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in J1 from "Valores_consumo Valor" to "Valores_consumo"
$ws.Range("J1").Value = "Valores_consumo"

# Move the active selection to H5 (matches the final selection state in the file)
$ws.Range("H5").Select()
